$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.502.07'
$ws.Range("E2").Value = '  +2.45%  '

$ws.Range("D3").Value = '3.477.06'
$ws.Range("E3").Value = '  +1.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.78'
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.43'
$ws.Range("E6").Value = '  +1.85%  '

$ws.Range("D7").Value = '3.476.56'

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("E9").Value = '  +0.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.69'
$ws.Range("E10").Value = '  +0.87%  '

$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("E12").Value = '  +4.67%  '

$ws.Range("D13").Value = '4.076.62'
$ws.Range("E13").Value = '  +1.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.66'
$ws.Range("E14").Value = '  +4.91%  '

$ws.Range("E15").Value = '  +2.49%  '

$ws.Range("D16").Value = '3.479.98'
$ws.Range("E16").Value = '  +1.73%  '

$ws.Range("E17").Value = '  +1.30%  '

$ws.Range("D18").Value = '63.566.05'
$ws.Range("E18").Value = '  +2.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.37'
$ws.Range("E19").Value = '  +2.83%  '

$ws.Range("E20").Value = '  +3.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.34'
$ws.Range("E21").Value = '  +1.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.21'
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("E23").Value = '  +2.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.91'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").Value = '3.621.03'
$ws.Range("E26").Value = '  +1.66%  '

$ws.Range("E27").Value = '  +0.54%  '

$ws.Range("E28").Value = '  -5.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.62'
$ws.Range("E29").Value = '  +2.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("E31").Value = '  +2.56%  '

$ws.Range("E32").Value = '  -0.56%  '

$ws.Range("E34").Value = '  -4.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.43'
$ws.Range("E35").Value = '  -0.52%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.17'
$ws.Range("E36").Value = '  +2.67%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.33'
$ws.Range("E37").Value = '  +1.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '31.68'
$ws.Range("E39").Value = '  +10.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '169.66'
$ws.Range("E40").Value = '  +0.76%  '

$ws.Range("D41").Value = '3.519.32'
$ws.Range("E41").Value = '  +1.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0765'
$ws.Range("E42").Value = '  +1.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.800'
$ws.Range("E43").Value = '  +1.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.41'
$ws.Range("E45").Value = '  -0.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  +3.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.42'
$ws.Range("E47").Value = '  -0.67%  '

$ws.Range("D48").Value = '2.610.21'
$ws.Range("E48").Value = '  +2.97%  '

$ws.Range("E49").Value = '  +10.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.12'
$ws.Range("E50").Value = '  +0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.80'
$ws.Range("E51").Value = '  +2.65%  '
